$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: move a single cell's value + formatting from $src to $dst, then
# blank out $src (value + style), mirroring a drag/cut-paste in the UI.
function Move-Cell($src, $dst) {
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $v = $ws.Range($src).Value()
    $ws.Range($dst).Value = $v
    $ws.Range($src).Clear() | Out-Null
}

# ---------------------------------------------------------------------
# 1) Shift the "coaches_table" and "gym_managers" ERD boxes (header +
#    the two helper cells in rows 19:20) one column to the right, to
#    make room for the new "gym_id" column being added to coaches_table.
# ---------------------------------------------------------------------
$ws.Range("I19:J20").UnMerge() | Out-Null
$ws.Range("L19:M20").UnMerge() | Out-Null

Move-Cell "M19" "N19"
Move-Cell "L19" "M19"
Move-Cell "M20" "N20"
Move-Cell "L20" "M20"

Move-Cell "J19" "K19"
Move-Cell "I19" "J19"
Move-Cell "J20" "K20"
Move-Cell "I20" "J20"

$ws.Range("J19:K20").Merge() | Out-Null
$ws.Range("M19:N20").Merge() | Out-Null

# ---------------------------------------------------------------------
# 2) Row 22 column headers: shift gym_managers' "manager_id"/"gym_id"
#    one column right, then add the new "gym_id" column header for
#    coaches_table (id/name stay put at I22/J22).
# ---------------------------------------------------------------------
Move-Cell "M22" "N22"
Move-Cell "L22" "M22"

$ws.Range("N22").Copy() | Out-Null
$ws.Range("K22").PasteSpecial(-4122) | Out-Null
$ws.Range("K22").Value = "gym_id"

# ---------------------------------------------------------------------
# 3) Remove the now-unused "admins" ERD box entirely (header box in
#    J24:K25 and its single "manager_id" column in J27:K27).
# ---------------------------------------------------------------------
$ws.Range("J24:K25").UnMerge() | Out-Null
$ws.Range("J24:K25").Clear() | Out-Null
$ws.Range("J27:K27").UnMerge() | Out-Null
$ws.Range("J27:K27").Clear() | Out-Null

# ---------------------------------------------------------------------
# 4) Cosmetic view state: scroll / selection left where the edit ended.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 13
$win.TopLeftCell = $ws.Range("E13")
$ws.Range("J23").Select() | Out-Null
